$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: insert the newest observation at row 232, pushing the
# existing historical rows (232:363) down to (233:364).
$ws.Rows("232:232").Insert()

$ws.Range("A232").Value = 5
$ws.Range("B232").Value = "Macroferia Regional de Talca"
$ws.Range("C232").Value = "Maule"
$ws.Range("D232").Value = 44873
$ws.Range("E232").Value = 7
$ws.Range("F232").Value = 100112008
$ws.Range("G232").Value = "Coliflor"
$ws.Range("H232").Value = "Sin especificar"
$ws.Range("I232").Value = "Primera"
$ws.Range("J232").Value = 3000
$ws.Range("K232").Value = 900
$ws.Range("L232").Value = 900
$ws.Range("M232").Value = 900
$ws.Range("N232").Value = "`$/unidad"
$ws.Range("O232").Value = "Región del Maule"
$ws.Range("P232").Value = 900
$ws.Range("Q232").Value = 1
$ws.Range("R232").Value = "Hortaliza"
